$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates that replace the raw date-serial values in column A (rows 2-27)
# with plain text values "2020-06-01" .. "2020-06-26" (shared strings,
# no cell style/number-format applied).
$dates = @(
    "2020-06-01","2020-06-02","2020-06-03","2020-06-04","2020-06-05",
    "2020-06-06","2020-06-07","2020-06-08","2020-06-09","2020-06-10",
    "2020-06-11","2020-06-12","2020-06-13","2020-06-14","2020-06-15",
    "2020-06-16","2020-06-17","2020-06-18","2020-06-19","2020-06-20",
    "2020-06-21","2020-06-22","2020-06-23","2020-06-24","2020-06-25",
    "2020-06-26"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    # Leading apostrophe forces literal-text entry so Excel doesn't
    # re-interpret the date-looking string back into a date serial.
    $ws.Cells.Item($row, 1).Value = "'" + $dates[$i]
    # Drop the old date number-format style so the cell goes back to
    # the workbook's default (unstyled) formatting.
    $ws.Cells.Item($row, 1).Style = "Normal"
}

# New row 28: raw + clean SSA data for June 27th.
$ws.Cells.Item(28, 1).Value = "'2020-06-27"
$ws.Cells.Item(28, 1).Style = "Normal"
$ws.Cells.Item(28, 2).Value = 212802
$ws.Cells.Item(28, 3).Value = 271151
$ws.Cells.Item(28, 4).Value = 67099
$ws.Cells.Item(28, 5).Value = 26381
$ws.Cells.Item(28, 6).Value = 31.24
